# Update odds values on Sheet1 to reflect latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value  = 1.85
$ws.Range("H2").Value  = 3.5
$ws.Range("I2").Value  = 4.2
$ws.Range("J2").Value  = 2.5
$ws.Range("Q2").Value  = 1.88
$ws.Range("R2").Value  = 1.98
$ws.Range("S2").Value  = 1.36
$ws.Range("T2").Value  = 3
$ws.Range("X2").Value  = 9
$ws.Range("Y2").Value  = 8.5
$ws.Range("Z2").Value  = 15
$ws.Range("AB2").Value = 23
$ws.Range("AC2").Value = 11
$ws.Range("AH2").Value = 13
$ws.Range("AJ2").Value = 15
$ws.Range("AL2").Value = 34
$ws.Range("AM2").Value = 41
$ws.Range("AO2").Value = 10
$ws.Range("AS2").Value = 126
$ws.Range("AT2").Value = 3
$ws.Range("AW2").Value = 6

# Row 3
$ws.Range("G3").Value = 4.1
$ws.Range("I3").Value = 2
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9

# Row 4
$ws.Range("Q4").Value = 2.2
$ws.Range("R4").Value = 1.65
